$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Update rows 2-4 with the new fixture data (game details + computed probabilities).
$ws.Range("A2").Value = 1369
$ws.Range("B2").Value = "2025-12-02T12:15:00"
$ws.Range("C2").Value = "Амур"
$ws.Range("D2").Value = "Динамо Мн"
$ws.Range("E2").Value = 897838
$ws.Range("F2").Value = "https://text.khl.ru/text/897838.html"
$ws.Range("G2").Value = 1.1
$ws.Range("H2").Value = 4.615385
$ws.Range("I2").Value = 1.552526
$ws.Range("J2").Value = 1.230769
$ws.Range("K2").Value = 1.165385
$ws.Range("L2").Value = 3.083955
$ws.Range("M2").Value = 5.715385
$ws.Range("N2").Value = 22.192716
$ws.Range("O2").Value = 42.630188
$ws.Range("P2").Value = 64.82290399999999
$ws.Range("Q2").Value = -0.2
$ws.Range("R2").Value = 0.2
$ws.Range("S2").Value = 0.110222
$ws.Range("T2").Value = 0.134865
$ws.Range("U2").Value = 0.7545460000000001
$ws.Range("V2").Value = 0.386332
$ws.Range("W2").Value = 0.613301
$ws.Range("X2").Value = 0.580246
$ws.Range("Y2").Value = 0.419386
$ws.Range("Z2").Value = 0.745048
$ws.Range("AA2").Value = 0.254585
$ws.Range("AB2").Value = 0.861764
$ws.Range("AC2").Value = 0.137868
$ws.Range("AD2").Value = 0.932617
$ws.Range("AE2").Value = 0.06701600000000001
$ws.Range("AF2").Value = 0.324827
$ws.Range("AG2").Value = 0.675173
$ws.Range("AH2").Value = 0.113094
$ws.Range("AI2").Value = 0.886906
$ws.Range("AJ2").Value = 0.813045
$ws.Range("AK2").Value = 0.186955
$ws.Range("AL2").Value = 0.595354
$ws.Range("AM2").Value = 0.404646
$ws.Range("AN2").Value = 0.432665
$ws.Range("AO2").Value = 0.960294
$ws.Range("A3").Value = 1369
$ws.Range("B3").Value = "2025-12-02T12:30:00"
$ws.Range("C3").Value = "Адмирал"
$ws.Range("D3").Value = "ХК Сочи"
$ws.Range("E3").Value = 897839
$ws.Range("F3").Value = "https://text.khl.ru/text/897839.html"
$ws.Range("G3").Value = 2.505511
$ws.Range("H3").Value = 0.928571
$ws.Range("I3").Value = 1.913907
$ws.Range("J3").Value = 6.535714
$ws.Range("K3").Value = 4.520612
$ws.Range("L3").Value = 1.421239
$ws.Range("M3").Value = 3.434082
$ws.Range("N3").Value = 32.125309
$ws.Range("O3").Value = 21.640231
$ws.Range("P3").Value = 53.76554
$ws.Range("Q3").Value = -0.038093
$ws.Range("R3").Value = -0.2
$ws.Range("S3").Value = 0.856225
$ws.Range("T3").Value = 0.07615
$ws.Range("U3").Value = 0.060739
$ws.Range("V3").Value = 0.156469
$ws.Range("W3").Value = 0.836645
$ws.Range("X3").Value = 0.292915
$ws.Range("Y3").Value = 0.700198
$ws.Range("Z3").Value = 0.455064
$ws.Range("AA3").Value = 0.538049
$ws.Range("AB3").Value = 0.615642
$ws.Range("AC3").Value = 0.377472
$ws.Range("AD3").Value = 0.751946
$ws.Range("AE3").Value = 0.241168
$ws.Range("AF3").Value = 0.939923
$ws.Range("AG3").Value = 0.060077
$ws.Range("AH3").Value = 0.828727
$ws.Range("AI3").Value = 0.171273
$ws.Range("AJ3").Value = 0.415477
$ws.Range("AK3").Value = 0.584523
$ws.Range("AL3").Value = 0.171658
$ws.Range("AM3").Value = 0.828342
$ws.Range("AN3").Value = 0.970588
$ws.Range("AO3").Value = 0.258435
$ws.Range("A4").Value = 1369
$ws.Range("B4").Value = "2025-12-02T19:00:00"
$ws.Range("C4").Value = "Локомотив"
$ws.Range("D4").Value = "СКА"
$ws.Range("E4").Value = 897840
$ws.Range("F4").Value = "https://text.khl.ru/text/897840.html"
$ws.Range("G4").Value = 2.392003
$ws.Range("H4").Value = 3.857143
$ws.Range("I4").Value = 2.506591
$ws.Range("J4").Value = 3.762853
$ws.Range("K4").Value = 3.077428
$ws.Range("L4").Value = 3.181867
$ws.Range("M4").Value = 6.249146
$ws.Range("N4").Value = 28.33909
$ws.Range("O4").Value = 33.415371
$ws.Range("P4").Value = 61.75446
$ws.Range("Q4").Value = -0.111154
$ws.Range("R4").Value = 0.2
$ws.Range("S4").Value = 0.401901
$ws.Range("T4").Value = 0.162865
$ws.Range("U4").Value = 0.434399
$ws.Range("V4").Value = 0.129522
$ws.Range("W4").Value = 0.869642
$ws.Range("X4").Value = 0.251846
$ws.Range("Y4").Value = 0.747318
$ws.Range("Z4").Value = 0.404979
$ws.Range("AA4").Value = 0.594185
$ws.Range("AB4").Value = 0.56473
$ws.Range("AC4").Value = 0.434434
$ws.Range("AD4").Value = 0.707577
$ws.Range("AE4").Value = 0.291588
$ws.Range("AF4").Value = 0.812122
$ws.Range("AG4").Value = 0.187878
$ws.Range("AH4").Value = 0.593931
$ws.Range("AI4").Value = 0.406069
$ws.Range("AJ4").Value = 0.826419
$ws.Range("AK4").Value = 0.173581
$ws.Range("AL4").Value = 0.616299
$ws.Range("AM4").Value = 0.383701
$ws.Range("AN4").Value = 0.716493
$ws.Range("AO4").Value = 0.744011

# Row 5's fixture was dropped from this tour's schedule; remove the row entirely
# so remaining data shifts up and the sheet dimension shrinks to A1:AO4.
$ws.Rows.Item(5).EntireRow.Delete()
